$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.415.96'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.849.84'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''240.66'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D8').Value = '''0.07655'
$ws.Range('E8').Value = '  +0.63%  '
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '''24.83'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').Value = '2.256.18'
$ws.Range('E11').Value = '  +22.51%  '
$ws.Range('D12').Value = '''0.07739'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '''5.042'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('E14').Value = '  +0.22%  '
$ws.Range('D15').Value = '''0.00001075'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '''83.36'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '''6.177'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '29.485.20'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').Value = '''228.27'
$ws.Range('E19').Value = '  -0.28%  '
$ws.Range('D20').Value = '''12.34'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '''7.469'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '''157.97'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').Value = '''0.1383'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('D26').Value = '''8.425'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').Value = '''1.386'
$ws.Range('E28').Value = '  +6.85%  '
$ws.Range('D29').Value = '''1.464'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').Value = '''4.134'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').Value = '''4.068'
$ws.Range('E32').Value = '  +0.92%  '
$ws.Range('D33').Value = '''1.842'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D35').Value = '''0.6967'
$ws.Range('E35').Value = '  -1.83%  '
$ws.Range('D36').Value = '''2.589'
$ws.Range('D37').Value = '''0.01805'
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').Value = '1.230.35'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').Value = '''2.724'
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').Value = '''6.454'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '''0.9092'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D43').Value = '''101.78'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '''66.08'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '''7.207'
$ws.Range('E45').Value = '  +0.52%  '
$ws.Range('D46').Value = '''0.00000000119'
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('D47').Value = '''0.4025'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '''8.993'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Value = '''1.683'
$ws.Range('E50').Value = '  +2.43%  '
$ws.Range('D51').Value = '''0.05702'
$ws.Range('E51').Value = '  -0.10%  '
